$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 17; existing rows 17..136 shift down to 18..137.
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the new observation.
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = 'Vega Monumental Concepción'
$ws.Range("C17").Value = 'Bíobío'
$ws.Range("D17").Value = 45168
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 100112012
$ws.Range("G17").Value = 'Espinaca'
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 7000
$ws.Range("L17").Value = 7000
$ws.Range("M17").Value = 7000
$ws.Range("N17").Value = '$/cuna 10 kilos'
$ws.Range("O17").Value = 'Región Metropolitana'
$ws.Range("P17").Value = 700
$ws.Range("Q17").Value = 10
$ws.Range("R17").Value = 'Hortaliza'
